$d = $word.ActiveDocument

# Two new centered, bold title paragraphs are inserted at the very top of the
# document, before the existing first paragraph ("Content"):
#   1) "Guidelines for Poster Submission "
#   2) "Geomundus 2019"
# Both use the same look: bold, color #333333, 14pt (sz/szCs 28), centered,
# with spacing-before of 220 (twips) and auto line rule - matching the style
# already used for the other headings in this document.

function Add-TitleParagraph([string]$text) {
    $target = $d.Paragraphs.Item(1).Range
    $target.Collapse(1)  # wdCollapseStart
    $target.InsertBefore($text + "`r")

    $newPara = $d.Paragraphs.Item(1).Range
    $newPara.Font.Bold = $true
    $newPara.Font.Color = 3355443   # 0x333333 -> RGB(51,51,51)
    $newPara.Font.Size = 14         # half-points 28 -> w:sz
    $newPara.Font.SizeBi = 14       # half-points 28 -> w:szCs
    $newPara.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter
}

Add-TitleParagraph "Geomundus 2019"
Add-TitleParagraph "Guidelines for Poster Submission "
